$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1069.3334
$ws.Range("I2").Value = 243.3077
$ws.Range("J2").Value = 3217
$ws.Range("K2").Value = 243.3077
$ws.Range("L2").Value = 3217
$ws.Range("M2").Value = -130.3077
$ws.Range("N2").Value = -3443
$ws.Range("H17").Value = 2535.1667
$ws.Range("I17").Value = 2495
$ws.Range("J17").Value = 2537.5293
$ws.Range("K17").Value = 7485
$ws.Range("L17").Value = 7612.5879
$ws.Range("M17").Value = -7317
$ws.Range("N17").Value = -7948.5879
$ws.Range("H21").Value = 14008.5
$ws.Range("I21").Value = 14008.5
$ws.Range("K21").Value = 14008.5
$ws.Range("M21").Value = -13540.5
$ws.Range("H23").Value = 14008.5
$ws.Range("I23").Value = 14008.5
$ws.Range("K23").Value = 14008.5
$ws.Range("M23").Value = -13774.5
$ws.Range("H32").Value = 2391.3333
$ws.Range("I32").Value = 1397
$ws.Range("K32").Value = 1397
$ws.Range("M32").Value = -1071
$ws.Range("H51").Value = 4499.3335
$ws.Range("I51").Value = 2799
$ws.Range("J51").Value = 5349.5
$ws.Range("K51").Value = 2799
$ws.Range("L51").Value = 5349.5
$ws.Range("M51").Value = -2315
$ws.Range("N51").Value = -6317.5
$ws.Range("H75").Value = 29990.75
$ws.Range("J75").Value = 29990.75
$ws.Range("L75").Value = 29990.75
$ws.Range("N75").Value = -31862.75
$ws.Range("H78").Value = 29990.75
$ws.Range("J78").Value = 29990.75
$ws.Range("L78").Value = 89972.25
$ws.Range("N78").Value = -99332.25
$ws.Range("H100").Value = 2150.1667
$ws.Range("I100").Value = 1860.2
$ws.Range("J100").Value = 3600
$ws.Range("K100").Value = 1860.2
$ws.Range("L100").Value = 3600
$ws.Range("M100").Value = -1319.2
$ws.Range("N100").Value = -4682
$ws.Range("H135").Value = 429.5
$ws.Range("I135").Value = 429.5
$ws.Range("K135").Value = 3865.5
$ws.Range("M135").Value = -1330.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 7014.5
$ws.Range("I33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("M33").Value = -4671
$ws.Range("H61").Value = 4912.25
$ws.Range("I61").Value = 2374.5
$ws.Range("J61").Value = 7450
$ws.Range("K61").Value = 2374.5
$ws.Range("L61").Value = 7450
$ws.Range("M61").Value = -2162.5
$ws.Range("N61").Value = -7874
$ws.Range("H74").Value = 1374.1428
$ws.Range("I74").Value = 1084
$ws.Range("J74").Value = 2099.5
$ws.Range("K74").Value = 1084
$ws.Range("L74").Value = 2099.5
$ws.Range("M74").Value = -210
$ws.Range("N74").Value = -3847.5
$ws.Range("H76").Value = 39821.75
$ws.Range("J76").Value = 39821.75
$ws.Range("L76").Value = 39821.75
$ws.Range("N76").Value = -40497.75
$ws.Range("H77").Value = 1374.1428
$ws.Range("I77").Value = 1084
$ws.Range("J77").Value = 2099.5
$ws.Range("K77").Value = 5420
$ws.Range("L77").Value = 10497.5
$ws.Range("M77").Value = -1052
$ws.Range("N77").Value = -19233.5
$ws.Range("H79").Value = 39821.75
$ws.Range("J79").Value = 39821.75
$ws.Range("L79").Value = 39821.75
$ws.Range("N79").Value = -42161.75
$ws.Range("H97").Value = 519.9286
$ws.Range("I97").Value = 553.25
$ws.Range("K97").Value = 553.25
$ws.Range("M97").Value = -57.25
$ws.Range("H122").Value = 1925.2
$ws.Range("I122").Value = 822.1429000000001
$ws.Range("J122").Value = 4499
$ws.Range("K122").Value = 2466.4287
$ws.Range("L122").Value = 13497
$ws.Range("M122").Value = -16.42870000000039
$ws.Range("N122").Value = -18397
$ws.Range("H132").Value = 1735.2222
$ws.Range("I132").Value = 1846
$ws.Range("J132").Value = 849
$ws.Range("K132").Value = 5538
$ws.Range("L132").Value = 2547
$ws.Range("M132").Value = -3008
$ws.Range("N132").Value = -7607
$ws.Range("H136").Value = 4912.25
$ws.Range("I136").Value = 2374.5
$ws.Range("J136").Value = 7450
$ws.Range("K136").Value = 7123.5
$ws.Range("L136").Value = 22350
$ws.Range("M136").Value = -4573.5
$ws.Range("N136").Value = -27450

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1569
$ws.Range("I99").Value = 1569
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1569
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -71
$ws.Range("N99").ClearContents()
$ws.Range("H130").Value = 91726.09
$ws.Range("J130").Value = 91726.09
$ws.Range("L130").Value = 91726.09
$ws.Range("N130").Value = -101766.09
$ws.Range("H134").Value = 815.6667
$ws.Range("I134").Value = 999
$ws.Range("J134").Value = 449
$ws.Range("K134").Value = 2997
$ws.Range("L134").Value = 1347
$ws.Range("M134").Value = -462
$ws.Range("N134").Value = -6417

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1674.9
$ws.Range("J31").Value = 1849.75
$ws.Range("L31").Value = 1849.75
$ws.Range("N31").Value = -2439.75
$ws.Range("H34").Value = 1674.9
$ws.Range("J34").Value = 1849.75
$ws.Range("L34").Value = 1849.75
$ws.Range("N34").Value = -2253.75
$ws.Range("H58").Value = 2683.5
$ws.Range("I58").Value = 2654.5
$ws.Range("K58").Value = 2654.5
$ws.Range("M58").Value = -2451.5
$ws.Range("H132").Value = 2175
$ws.Range("I132").Value = 2096.875
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 6290.625
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -3760.625
$ws.Range("N132").Value = -13460
$ws.Range("H133").Value = 44998.25
$ws.Range("H134").Value = 1093.4445
$ws.Range("I134").Value = 1093.4445
$ws.Range("K134").Value = 3280.3335
$ws.Range("M134").Value = -745.3335000000002
$ws.Range("H136").Value = 2683.5
$ws.Range("I136").Value = 2654.5
$ws.Range("K136").Value = 7963.5
$ws.Range("M136").Value = -5413.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 59983.668
$ws.Range("J37").Value = 59983.668
$ws.Range("L37").Value = 179951.004
$ws.Range("N37").Value = -180175.004
$ws.Range("H113").Value = 1575.3077
$ws.Range("I113").Value = 1188.4
$ws.Range("J113").Value = 1817.125
$ws.Range("K113").Value = 3565.2
$ws.Range("L113").Value = 5451.375
$ws.Range("M113").Value = -1395.2
$ws.Range("N113").Value = -9791.375
$ws.Range("H121").Value = 16091.909
$ws.Range("J121").Value = 6331.222
$ws.Range("L121").Value = 18993.666
$ws.Range("N121").Value = -21613.666

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 58032
$ws.Range("I62").Value = 58032
$ws.Range("K62").Value = 58032
$ws.Range("M62").Value = -57346
$ws.Range("H65").Value = 58032
$ws.Range("I65").Value = 58032
$ws.Range("K65").Value = 174096
$ws.Range("M65").Value = -170664
$ws.Range("H122").Value = 3312.5715
$ws.Range("J122").Value = 2500
$ws.Range("L122").Value = 7500
$ws.Range("N122").Value = -12400
$ws.Range("H126").Value = 3799
$ws.Range("I126").Value = 3598.3333
$ws.Range("J126").Value = 4100
$ws.Range("K126").Value = 10794.9999
$ws.Range("L126").Value = 12300
$ws.Range("M126").Value = -8324.999899999999
$ws.Range("N126").Value = -17240

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 88998
$ws.Range("J63").Value = 88998
$ws.Range("L63").Value = 88998
$ws.Range("N63").Value = -90496
$ws.Range("H66").Value = 88998
$ws.Range("J66").Value = 88998
$ws.Range("L66").Value = 266994
$ws.Range("N66").Value = -274482
$ws.Range("H122").Value = 4489.1304
$ws.Range("I122").Value = 3522.4375
$ws.Range("K122").Value = 10567.3125
$ws.Range("M122").Value = -8117.3125
$ws.Range("H136").Value = 1839.8
$ws.Range("I136").Value = 799.75
$ws.Range("K136").Value = 2399.25
$ws.Range("M136").Value = 150.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 44720.5
$ws.Range("J54").Value = 44720.5
$ws.Range("L54").Value = 44720.5
$ws.Range("N54").Value = -45760.5
$ws.Range("H75").Value = 90130
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 90130
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 90130
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -92002
$ws.Range("H78").Value = 90130
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 90130
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 270390
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -279750
$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -21490
$ws.Range("H113").Value = 1285.625
$ws.Range("I113").Value = 1224
$ws.Range("J113").Value = 1470.5
$ws.Range("K113").Value = 3672
$ws.Range("L113").Value = 4411.5
$ws.Range("M113").Value = -1502
$ws.Range("N113").Value = -8751.5
$ws.Range("H122").Value = 983
$ws.Range("I122").Value = 983
$ws.Range("K122").Value = 2949
$ws.Range("M122").Value = -499
$ws.Range("H126").Value = 2336
$ws.Range("I126").Value = 2261.4167
$ws.Range("K126").Value = 6784.250100000001
$ws.Range("M126").Value = -4314.250100000001
$ws.Range("H132").Value = 5224.75
$ws.Range("I132").Value = 5436.091
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 16308.273
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -13778.273
$ws.Range("N132").Value = -13760

Write-Host "Applied all profit updates"